$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update / add cell values
$ws.Range("D2").Value = -0.0509
$ws.Range("E2").Value = -0.134
$ws.Range("G2").Value = 0.05091278399117308
$ws.Range("H2").Value = 0.05091278399117308
$ws.Range("I2").Value = 0.005958172425898991
$ws.Range("J2").Value = 0.005556456391203659
$ws.Range("K2").Value = -39.42
$ws.Range("L2").Value = -0.01977029941321029
$ws.Range("M2").Value = 23.4
$ws.Range("N2").Value = 0.012464046021093
$ws.Range("O2").Value = -0.593607305936073
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 23.4
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 261.49
$ws.Range("V2").Value = 0.1392830510280175
$ws.Range("W2").Value = -0.06989525029750179
$ws.Range("X2").Value = 0.06745754649294133
$ws.Range("Y2").Value = -0.1373527967904431
$ws.Range("Z2").Value = 0.9481938711456889
$ws.Range("AA2").Value = -0.009202690449714345
$ws.Range("AB2").Value = 0.05985323818331078
$ws.Range("AC2").Value = -0.06948120837974261
$ws.Range("AD2").Value = 1137.361
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1137.361
$ws.Range("AG2").Value = 875.8710000000001
$ws.Range("AH2").Value = 0.3772640683623014
$ws.Range("AI2").Value = 0.5262244483915459
$ws.Range("AJ2").Value = 0.3181201559890036
$ws.Range("AK2").Value = 0.4610160374046449
$ws.Range("AL2").Value = 32.09
$ws.Range("AM2").Value = 32.09
$ws.Range("AN2").Value = 14.35880570635021
$ws.Range("AO2").Value = 0.3702087877843565
$ws.Range("AP2").Value = 11.05758111349577
$ws.Range("AQ2").Value = 0.3702087877843565
$ws.Range("D3").Value = -0.0228
$ws.Range("E3").Value = -0.064
$ws.Range("G3").Value = 0.05765595463137996
$ws.Range("H3").Value = 0.05765595463137996
$ws.Range("I3").Value = 0.08695652173913043
$ws.Range("J3").Value = 0.07067868504772004
$ws.Range("K3").Value = 13.3
$ws.Range("L3").Value = 0.06285444234404537
$ws.Range("M3").Value = 23.4
$ws.Range("N3").Value = 0.07186732186732185
$ws.Range("O3").Value = 1.759398496240601
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 23.4
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 66.59999999999999
$ws.Range("V3").Value = 0.2045454545454545
$ws.Range("W3").Value = 0.09729334308705195
$ws.Range("X3").Value = 0.05982182809544657
$ws.Range("Y3").Value = 0.03747151499160538
$ws.Range("Z3").Value = 3.009957325746799
$ws.Range("AA3").Value = 0.2127398258335357
$ws.Range("AB3").Value = 0.05970562628505031
$ws.Range("AC3").Value = 0.1530341995484854
$ws.Range("AD3").Value = 0.961
$ws.Range("AF3").Value = 0.961
$ws.Range("AG3").Value = -65.639
$ws.Range("AH3").Value = 0.002942788636732494
$ws.Range("AI3").Value = 0.008105532173311629
$ws.Range("AJ3").Value = -0.2524955666426887
$ws.Range("AK3").Value = -1.263235888454803
$ws.Range("AN3").Value = 0.04368181818181818
$ws.Range("AP3").Value = -2.983590909090909
$ws.Range("D4").Value = -0.0156
$ws.Range("E4").Value = -0.204
$ws.Range("G4").Value = 0.08274753966659973
$ws.Range("H4").Value = 0.08274753966659973
$ws.Range("I4").Value = 0.05985137577826873
$ws.Range("J4").Value = 0.04684321549351206
$ws.Range("K4").Value = 13.5
$ws.Range("L4").Value = 0.02711387828881302
$ws.Range("O4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0.08549715009499682
$ws.Range("X4").Value = 0.05970281800183537
$ws.Range("Y4").Value = 0.02579433209316145
$ws.Range("Z4").Value = 1.477272727272727
$ws.Range("AA4").Value = 0.06920020470632463
$ws.Range("AB4").Value = 0.05970281800183537
$ws.Range("AC4").Value = 0.009497386704489255
$ws.Range("AD4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 0
$ws.Range("AK4").Value = 0
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AO4").Value = 0
$ws.Range("AQ4").Value = 0
$ws.Range("D5").Value = 0.0225
$ws.Range("G5").Value = 0.2149046793760832
$ws.Range("H5").Value = 0.2149046793760832
$ws.Range("I5").Value = -0.009358752166377816
$ws.Range("J5").Value = -0.009358752166377816
$ws.Range("K5").Value = -9.550000000000001
$ws.Range("L5").Value = -0.08275563258232235
$ws.Range("U5").Value = 2.89
$ws.Range("V5").Value = 0.01510716152639833
$ws.Range("W5").Value = -0.06831187410586552
$ws.Range("X5").Value = 0.09924518562188478
$ws.Range("Y5").Value = -0.1675570597277503
$ws.Range("Z5").Value = 0.3812355467459531
$ws.Range("AA5").Value = -0.00356788899900892
$ws.Range("AB5").Value = 0.06017530521805391
$ws.Range("AC5").Value = -0.06374319421706283
$ws.Range("AD5").Value = 187.6
$ws.Range("AF5").Value = 187.6
$ws.Range("AG5").Value = 184.71
$ws.Range("AH5").Value = 0.4951174452362101
$ws.Range("AI5").Value = 0.590122680088078
$ws.Range("AJ5").Value = 0.491236935187894
$ws.Range("AK5").Value = 0.5863623377035651
$ws.Range("AL5").Value = 8.210000000000001
$ws.Range("AM5").Value = 8.210000000000001
$ws.Range("AN5").Value = 28.77300613496933
$ws.Range("AO5").Value = -0.1315468940316687
$ws.Range("AP5").Value = 28.329754601227
$ws.Range("AQ5").Value = -0.1315468940316687
$ws.Range("D6").Value = -0.079
$ws.Range("G6").Value = 0.0203997949769349
$ws.Range("H6").Value = 0.0203997949769349
$ws.Range("I6").Value = -0.01619682214249103
$ws.Range("J6").Value = -0.01619682214249103
$ws.Range("K6").Value = -37.2
$ws.Range("L6").Value = -0.03813429010763711
$ws.Range("U6").Value = 167.6
$ws.Range("V6").Value = 0.5048192771084338
$ws.Range("W6").Value = -0.1080453093232646
$ws.Range("X6").Value = 0.1591117494102745
$ws.Range("Y6").Value = -0.2671570587335391
$ws.Range("Z6").Value = 0.9160742625860434
$ws.Range("AA6").Value = -0.01483749190041977
$ws.Range("AB6").Value = 0.06038173064200262
$ws.Range("AC6").Value = -0.07521922254242239
$ws.Range("AD6").Value = 818.5
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 818.5
$ws.Range("AG6").Value = 650.9
$ws.Range("AH6").Value = 0.7114298131247284
$ws.Range("AI6").Value = 0.6948217317487266
$ws.Range("AJ6").Value = 0.6622240309288839
$ws.Range("AK6").Value = 0.644200316706255
$ws.Range("AL6").Value = 19.6
$ws.Range("AM6").Value = 19.6
$ws.Range("AN6").Value = 63.44961240310077
$ws.Range("AO6").Value = -0.8061224489795918
$ws.Range("AP6").Value = 50.45736434108527
$ws.Range("AQ6").Value = -0.8061224489795918
$ws.Range("D7").Value = -0.0978
$ws.Range("G7").Value = 0.0207380073800738
$ws.Range("H7").Value = 0.0207380073800738
$ws.Range("I7").Value = -0.0929889298892989
$ws.Range("J7").Value = -0.0929889298892989
$ws.Range("K7").Value = -10.2
$ws.Range("L7").Value = -0.07527675276752767
$ws.Range("U7").Value = 11.1
$ws.Range("V7").Value = 0.02958422174840085
$ws.Range("W7").Value = -0.07147862648913805
$ws.Range("X7").Value = 0.07063238429349135
$ws.Range("Y7").Value = -0.1421110107826294
$ws.Range("Z7").Value = 0.6138165345413363
$ws.Range("AA7").Value = -0.05707814269535674
$ws.Range("AB7").Value = 0.05990632318001283
$ws.Range("AC7").Value = -0.1169844658753696
$ws.Range("AD7").Value = 101.7
$ws.Range("AF7").Value = 101.7
$ws.Range("AG7").Value = 90.60000000000001
$ws.Range("AH7").Value = 0.2132522541413294
$ws.Range("AI7").Value = 0.4357326478149101
$ws.Range("AJ7").Value = 0.1945040790038643
$ws.Range("AK7").Value = 0.407557354925776
$ws.Range("AL7").Value = 3.28
$ws.Range("AM7").Value = 3.28
$ws.Range("AN7").Value = -97.78846153846153
$ws.Range("AO7").Value = -3.841463414634147
$ws.Range("AP7").Value = -87.11538461538463
$ws.Range("AQ7").Value = -3.841463414634147
$ws.Range("D8").Value = -0.122
$ws.Range("G8").Value = 0.01043103448275862
$ws.Range("H8").Value = 0.01043103448275862
$ws.Range("I8").Value = -0.1179310344827586
$ws.Range("J8").Value = -0.1179310344827586
$ws.Range("K8").Value = -9.27
$ws.Range("L8").Value = -0.1598275862068965
$ws.Range("M8").Value = -0
$ws.Range("N8").Value = -0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = -0
$ws.Range("Q8").Value = -0
$ws.Range("R8").Value = 0
$ws.Range("U8").Value = 13.3
$ws.Range("V8").Value = 0.0528196981731533
$ws.Range("W8").Value = -0.1081680280046674
$ws.Range("X8").Value = 0.06428270869239133
$ws.Range("Y8").Value = -0.1724507366970588
$ws.Range("Z8").Value = 0.5411457361448031
$ws.Range("AA8").Value = -0.06381787646949057
$ws.Range("AB8").Value = 0.05980015318660873
$ws.Range("AC8").Value = -0.1236180296560993
$ws.Range("AD8").Value = 28.6
$ws.Range("AF8").Value = 28.6
$ws.Range("AG8").Value = 15.3
$ws.Range("AH8").Value = 0.1019971469329529
$ws.Range("AI8").Value = 0.2718631178707225
$ws.Range("AJ8").Value = 0.05728191688506177
$ws.Range("AK8").Value = 0.1664853101196953
$ws.Range("AL8").Value = 1
$ws.Range("AM8").Value = 1
$ws.Range("AN8").Value = -12.06751054852321
$ws.Range("AO8").Value = -6.84
$ws.Range("AP8").Value = -6.455696202531645
$ws.Range("AQ8").Value = -6.84

# Remove cells that no longer exist in the target
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("T8").ClearContents()
